# Insert a new row at position 115 on the active sheet, shifting existing
# rows 115-141 down to 116-142, then populate the newly inserted row with
# the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 115 (shifts rows 115:141 -> 116:142)
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with the new record's data
$ws.Cells.Item(115, 1).Value = 9
$ws.Cells.Item(115, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(115, 3).Value = "Metropolitana"
$ws.Cells.Item(115, 4).Value = 44551
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100101
$ws.Cells.Item(115, 8).Value = "Berries"
$ws.Cells.Item(115, 9).Value = 100101001
$ws.Cells.Item(115, 10).Value = "Arándano (blue)"
$ws.Cells.Item(115, 11).Value = "Sin especificar"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 250
$ws.Cells.Item(115, 14).Value = 3000
$ws.Cells.Item(115, 15).Value = 3000
$ws.Cells.Item(115, 16).Value = 3000
$ws.Cells.Item(115, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(115, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(115, 19).Value = 1500
$ws.Cells.Item(115, 20).Value = 2
